# Refresh the "Authors" (column E) values for rows 2-22 (except row 10,
# which was already correct) with their corrected/re-cleaned text.
# The content (names/emails) is unchanged; this mirrors the upstream
# re-clean pass that re-wrote these entries (commit: "Cleaning more data.")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = '[Barnaby Edward%Young%xref no email%0,        Sean Wei Xiang%Ong%xref no email%0,        Shirin%Kalimuddin%xref no email%0,        Jenny G.%Low%xref no email%0,        Seow Yen%Tan%xref no email%0,        Jiashen%Loh%xref no email%0,        Oon-Tek%Ng%xref no email%0,        Kalisvar%Marimuthu%xref no email%0,        Li Wei%Ang%xref no email%0,        Tze Minn%Mak%xref no email%0,        Sok Kiang%Lau%xref no email%0,        Danielle E.%Anderson%xref no email%0,        Kian Sing%Chan%xref no email%0,        Thean Yen%Tan%xref no email%0,        Tong Yong%Ng%xref no email%0,        Lin%Cui%xref no email%0,        Zubaidah%Said%xref no email%0,        Lalitha%Kurupatham%xref no email%0,        Mark I-Cheng%Chen%xref no email%0,        Monica%Chan%xref no email%0,        Shawn%Vasoo%xref no email%0,        Lin-Fa%Wang%xref no email%0,        Boon Huan%Tan%xref no email%0,        Raymond Tzer Pin%Lin%xref no email%0,        Vernon Jian Ming%Lee%xref no email%0,        Yee-Sin%Leo%xref no email%0,        David Chien%Lye%xref no email%0]'
$ws.Range("E3").Value = '[Fei%Zhou%NULL%0,                Ting%Yu%NULL%0,                Ronghui%Du%NULL%0,                Guohui%Fan%NULL%0,                Ying%Liu%NULL%0,                Zhibo%Liu%NULL%0,                Jie%Xiang%NULL%0,                Yeming%Wang%NULL%0,                Bin%Song%NULL%0,                Xiaoying%Gu%NULL%0,                Lulu%Guan%NULL%0,                Yuan%Wei%NULL%0,                Hui%Li%NULL%0,                Xudong%Wu%NULL%0,                Jiuyang%Xu%NULL%0,                Shengjin%Tu%NULL%0,                Yi%Zhang%NULL%0,                Hua%Chen%NULL%0,                Bin%Cao%NULL%0]'
$ws.Range("E4").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,               Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%0,               Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%0,               Subert-Salas%Lizandra%coreGivesNoEmail%0]'
$ws.Range("E5").Value = '[Ning%Tang%NULL%0,                Dengju%Li%NULL%0,                Xiong%Wang%NULL%0,                Ziyong%Sun%NULL%0]'
$ws.Range("E6").Value = '[Rong%Qu%NULL%1,                Yun%Ling%NULL%0,                Yi‐hui‐zhi%Zhang%NULL%1,                Li‐ya%Wei%NULL%1,                Xiao%Chen%NULL%1,                Xu‐mian%Li%NULL%1,                Xuan‐yong%Liu%NULL%1,                Han‐mian%Liu%NULL%1,                Zhi%Guo%guozhi77@126.com%1,                Hua%Ren%renhua2009@hotmail.com%2,                Hua%Ren%renhua2009@hotmail.com%0,                Qiang%Wang%wangqiang@wust.edu.cn%0]'
$ws.Range("E7").Value = '[Qiurong%Ruan%NULL%0,                Kun%Yang%NULL%0,                Kun%Yang%NULL%0,                Wenxia%Wang%NULL%0,                Wenxia%Wang%NULL%0,                Lingyu%Jiang%NULL%0,                Lingyu%Jiang%NULL%0,                Jianxin%Song%songsingsjx@sina.com%0,                Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E8").Value = '[Chaolin%Huang%NULL%0,                Yeming%Wang%NULL%0,                Xingwang%Li%NULL%0,                Lili%Ren%NULL%0,                Jianping%Zhao%NULL%0,                Yi%Hu%NULL%0,                Li%Zhang%NULL%0,                Guohui%Fan%NULL%0,                Jiuyang%Xu%NULL%0,                Xiaoying%Gu%NULL%0,                Zhenshun%Cheng%NULL%0,                Ting%Yu%NULL%0,                Jiaan%Xia%NULL%0,                Yuan%Wei%NULL%0,                Wenjuan%Wu%NULL%0,                Xuelei%Xie%NULL%0,                Wen%Yin%NULL%0,                Hui%Li%NULL%0,                Min%Liu%NULL%0,                Yan%Xiao%NULL%0,                Hong%Gao%NULL%0,                Li%Guo%NULL%0,                Jungang%Xie%NULL%0,                Guangfa%Wang%NULL%0,                Rongmeng%Jiang%NULL%0,                Zhancheng%Gao%NULL%0,                Qi%Jin%NULL%0,                Jianwei%Wang%wangjw28@163.com%0,                Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("E9").Value = '[Yingxia%Liu%yingxialiu@hotmail.com%0,                Yang%Yang%NULL%0,                Cong%Zhang%NULL%2,                Fengming%Huang%NULL%2,                Fuxiang%Wang%NULL%2,                Jing%Yuan%NULL%2,                Zhaoqin%Wang%NULL%2,                Jinxiu%Li%NULL%2,                Jianming%Li%NULL%2,                Cheng%Feng%NULL%2,                Zheng%Zhang%NULL%2,                Lifei%Wang%NULL%2,                Ling%Peng%NULL%2,                Li%Chen%NULL%2,                Yuhao%Qin%NULL%2,                Dandan%Zhao%NULL%2,                Shuguang%Tan%NULL%2,                Lu%Yin%NULL%3,                Jun%Xu%NULL%2,                Congzhao%Zhou%NULL%2,                Chengyu%Jiang%jiang@pumc.edu.cn%2,                Lei%Liu%liulei3322@aliyun.com%0]'
$ws.Range("E11").Value = '[Yong%Gao%NULL%0,                Tuantuan%Li%NULL%0,                Mingfeng%Han%NULL%0,                Xiuyong%Li%NULL%0,                Dong%Wu%NULL%0,                Yuanhong%Xu%NULL%0,                Yulin%Zhu%NULL%0,                Yan%Liu%NULL%0,                Xiaowu%Wang%wangxiaowu19880218@126.com%0,                Linding%Wang%wanglinding@ahmu.edu.cn%0]'
$ws.Range("E12").Value = '[Wei-jie%Guan%NULL%0,                Zheng-yi%Ni%NULL%13,                Zheng-yi%Ni%NULL%0,                Yu%Hu%NULL%0,                Wen-hua%Liang%NULL%6,                Chun-quan%Ou%NULL%6,                Jian-xing%He%NULL%8,                Lei%Liu%NULL%0,                Hong%Shan%NULL%0,                Chun-liang%Lei%NULL%8,                David S.C.%Hui%NULL%5,                Bin%Du%NULL%5,                Lan-juan%Li%NULL%5,                Guang%Zeng%NULL%5,                Kwok-Yung%Yuen%NULL%0,                Ru-chong%Chen%NULL%8,                Chun-li%Tang%NULL%8,                Tao%Wang%NULL%0,                Ping-yan%Chen%NULL%6,                Jie%Xiang%NULL%0,                Shi-yue%Li%NULL%8,                Jin-lin%Wang%NULL%5,                Zi-jing%Liang%NULL%5,                Yi-xiang%Peng%NULL%8,                Li%Wei%NULL%0,                Yong%Liu%NULL%0,                Ya-hua%Hu%NULL%8,                Peng%Peng%NULL%0,                Jian-ming%Wang%NULL%8,                Ji-yang%Liu%NULL%8,                Zhong%Chen%NULL%0,                Gang%Li%NULL%0,                Zhi-jian%Zheng%NULL%8,                Shao-qin%Qiu%NULL%8,                Jie%Luo%NULL%0,                Chang-jiang%Ye%NULL%8,                Shao-yong%Zhu%NULL%8,                Nan-shan%Zhong%NULL%8]'
$ws.Range("E13").Value = '[Huan%Han%xref no email%3,        Lan%Yang%xref no email%1,        Rui%Liu%xref no email%1,        Fang%Liu%xref no email%1,        Kai-lang%Wu%xref no email%1,        Jie%Li%xref no email%1,        Xing-hui%Liu%xref no email%1,        Cheng-liang%Zhu%xref no email%1]'
$ws.Range("E14").Value = '[Kunhua%Li%NULL%0,                Jiong%Wu%NULL%0,                Faqi%Wu%NULL%0,                Dajing%Guo%NULL%0,                Linli%Chen%NULL%0,                Zheng%Fang%NULL%0,                Chuanming%Li%NULL%0]'
$ws.Range("E15").Value = '[Wei%Liu%NULL%0,                Zhao-Wu%Tao%NULL%0,                Lei%Wang%NULL%0,                Ming-Li%Yuan%NULL%0,                Kui%Liu%NULL%0,                Ling%Zhou%NULL%0,                Shuang%Wei%NULL%0,                Yan%Deng%NULL%0,                Jing%Liu%NULL%0,                Hui-Guo%Liu%NULL%0,                Ming%Yang%NULL%0,                Yi%Hu%NULL%0,                Pei-Fang%Wei%NULL%0,                Pei-Fang%Wei%NULL%0]'
$ws.Range("E16").Value = '[Chuan%Qin%NULL%0,                Luoqi%Zhou%NULL%1,                Ziwei%Hu%NULL%1,                Shuoqi%Zhang%NULL%1,                Sheng%Yang%NULL%1,                Yu%Tao%NULL%0,                Cuihong%Xie%NULL%1,                Ke%Ma%NULL%0,                Ke%Shang%NULL%1,                Wei%Wang%NULL%0,                Dai-Shi%Tian%tiands@tjh.tjmu.edu.cn%1]'
$ws.Range("E17").Value = '[Dan%Sun%NULL%0,                Hui%Li%NULL%0,                Xiao-Xia%Lu%NULL%1,                Han%Xiao%NULL%1,                Jie%Ren%NULL%1,                Fu-Rong%Zhang%792523496@qq.com%1,                Zhi-Sheng%Liu%liuzsc@126.com%1]'
$ws.Range("E18").Value = '[Suxin%Wan%NULL%0,                Yi%Xiang%NULL%3,                Wei%Fang%NULL%3,                Yu%Zheng%yuzheng1@cdutcm.edu.cn%3,                Boqun%Li%279685211@qq.com%6,                Boqun%Li%279685211@qq.com%0,                Yanjun%Hu%huyanjun@163.com%3,                Chunhui%Lang%NULL%4,                Daoqiu%Huang%NULL%3,                Qiuyan%Sun%NULL%3,                Yan%Xiong%NULL%0,                Xia%Huang%NULL%3,                Jinglong%Lv%NULL%4,                Yaling%Luo%NULL%3,                Li%Shen%NULL%3,                Haoran%Yang%NULL%3,                Gu%Huang%NULL%3,                Ruishan%Yang%NULL%3]'
$ws.Range("E19").Value = '[Zhongliang%Wang%NULL%0,                Bohan%Yang%NULL%0,                Qianwen%Li%NULL%0,                Lu%Wen%NULL%0,                Ruiguang%Zhang%zrg27@163.com%0]'
$ws.Range("E20").Value = '[Chaomin%Wu%xref no email%0,        Xiaoyan%Chen%xref no email%0,        Yanping%Cai%xref no email%0,        Jia\u2019an%Xia%xref no email%0,        Xing%Zhou%xref no email%0,        Sha%Xu%xref no email%0,        Hanping%Huang%xref no email%0,        Li%Zhang%xref no email%0,        Xia%Zhou%xref no email%0,        Chunling%Du%xref no email%0,        Yuye%Zhang%xref no email%0,        Juan%Song%xref no email%0,        Sijiao%Wang%xref no email%0,        Yencheng%Chao%xref no email%0,        Zeyong%Yang%xref no email%0,        Jie%Xu%xref no email%0,        Xin%Zhou%xref no email%0,        Dechang%Chen%xref no email%0,        Weining%Xiong%xref no email%0,        Lei%Xu%xref no email%0,        Feng%Zhou%xref no email%0,        Jinjun%Jiang%xref no email%0,        Chunxue%Bai%xref no email%0,        Junhua%Zheng%xref no email%0,        Yuanlin%Song%xref no email%0]'
$ws.Range("E21").Value = '[Xiaobo%Yang%NULL%0,                Yuan%Yu%NULL%0,                Jiqian%Xu%NULL%0,                Huaqing%Shu%NULL%0,                Jia''an%Xia%NULL%0,                Hong%Liu%NULL%0,                Yongran%Wu%NULL%0,                Lu%Zhang%NULL%0,                Zhui%Yu%NULL%0,                Minghao%Fang%NULL%0,                Ting%Yu%NULL%0,                Yaxin%Wang%NULL%0,                Shangwen%Pan%NULL%0,                Xiaojing%Zou%NULL%0,                Shiying%Yuan%NULL%0,                You%Shang%NULL%0]'
$ws.Range("E22").Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,               Cao%Yi-yuan%coreGivesNoEmail%2,               Dong%Xiang%coreGivesNoEmail%2,               Gao%Ya-dong%coreGivesNoEmail%2,               Yan%You-qin%coreGivesNoEmail%2,               Yang%Yi-bin%coreGivesNoEmail%2,               Yuan%Ya-dong%coreGivesNoEmail%2,               Zhang%Jin-jin%coreGivesNoEmail%2]'
